$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Corrections to existing rows 138-143 ---
$ws.Range("P138").Value = 1807779
$ws.Range("AA138").Value = -28704
$ws.Range("AB138").Value = 473232

$ws.Range("J139").Value = -81642
$ws.Range("P139").Value = -261281
$ws.Range("AA139").Value = -25104
$ws.Range("AB139").Value = -324827

$ws.Range("B140").Value = -2406542
$ws.Range("J140").Value = -120340
$ws.Range("P140").Value = -750378
$ws.Range("AA140").Value = -25669
$ws.Range("AB140").Value = -1656164

$ws.Range("B141").Value = 2661805
$ws.Range("J141").Value = 2637872
$ws.Range("P141").Value = 1250785
$ws.Range("AA141").Value = -25019
$ws.Range("AB141").Value = 1411020

$ws.Range("P142").Value = 2165766
$ws.Range("AA142").Value = -25874
$ws.Range("AB142").Value = -3347299

$ws.Range("B143").Value = -2152217
$ws.Range("J143").Value = 41685
$ws.Range("P143").Value = 105869
$ws.Range("AA143").Value = -26447
$ws.Range("AB143").Value = -2258086

# --- New row 144 ---
# Leading apostrophe forces literal-text entry (column A stores its dates
# as plain text, e.g. "01-06-2021", not as Excel date serials). Reset the
# style back to Normal afterward so the cell has no explicit style index,
# matching its neighbours in column A.
$ws.Range("A144").Value = "'01-07-2021"
$ws.Range("A144").Style = "Normal"
$ws.Range("B144").Value = 2744853
$ws.Range("C144").Value = 125601
$ws.Range("D144").Value = 203581
$ws.Range("E144").Value = 77980
$ws.Range("F144").Value = 2441164
$ws.Range("G144").Value = 3302075
$ws.Range("H144").Value = 860911
$ws.Range("I144").Value = 30477
$ws.Range("J144").Value = 147611
$ws.Range("K144").Value = 0
$ws.Range("L144").Value = 0
$ws.Range("M144").Value = 0
$ws.Range("N144").Value = 0
$ws.Range("O144").Value = 0
$ws.Range("P144").Value = 6179213
$ws.Range("Q144").Value = 4325590
$ws.Range("R144").Value = 4327211
$ws.Range("S144").Value = 4327211
$ws.Range("T144").Value = 0
$ws.Range("U144").Value = 1622
$ws.Range("V144").Value = 1881263
$ws.Range("W144").Value = 2794447
$ws.Range("X144").Value = 2794447
$ws.Range("Y144").Value = 0
$ws.Range("Z144").Value = 913184
$ws.Range("AA144").Value = -27639
$ws.Range("AB144").Value = -3434360
